$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "In Translation" -> "Handed back: in sync with en-US"
#    Appears on all three sheets (Overview E/F cols, zh-cn + de-de "Status"
#    column C) for both data rows.
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: the handback finished, so "Latest Target File" (I) and
#    "Latest Handback File" (J) get filled in with links/filenames, and the
#    "Latest Handback DateTime" (K) is stamped.
# ---------------------------------------------------------------------------
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/20902bde2b74d45631a5c6772b2533a6aabf985e/e2e/6078e45d-6b8d-4f2a-b087-85fa4331d5cc.md",
    [Type]::Missing,
    [Type]::Missing,
    "6078e45d-6b8d-4f2a-b087-85fa4331d5cc.md"
)
$wsZhCn.Range("J2").Value = "6078e45d-6b8d-4f2a-b087-85fa4331d5cc.7e6d32f1bb119aec7d9a809334d3fae3a6c49a2b.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-28 12:24:50"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/20902bde2b74d45631a5c6772b2533a6aabf985e/e2e/7e314904-3a85-4fe1-a353-d521f0df730f.md",
    [Type]::Missing,
    [Type]::Missing,
    "7e314904-3a85-4fe1-a353-d521f0df730f.md"
)
$wsZhCn.Range("J3").Value = "7e314904-3a85-4fe1-a353-d521f0df730f.a47c97363969567fabfb364e0821a2933e6f6093.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-28 12:24:50"

# ---------------------------------------------------------------------------
# 3. de-de sheet: same handback completion.
# ---------------------------------------------------------------------------
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/20902bde2b74d45631a5c6772b2533a6aabf985e/e2e/6078e45d-6b8d-4f2a-b087-85fa4331d5cc.md",
    [Type]::Missing,
    [Type]::Missing,
    "6078e45d-6b8d-4f2a-b087-85fa4331d5cc.md"
)
$wsDeDe.Range("J2").Value = "6078e45d-6b8d-4f2a-b087-85fa4331d5cc.7e6d32f1bb119aec7d9a809334d3fae3a6c49a2b.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-28 12:24:57"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/20902bde2b74d45631a5c6772b2533a6aabf985e/e2e/7e314904-3a85-4fe1-a353-d521f0df730f.md",
    [Type]::Missing,
    [Type]::Missing,
    "7e314904-3a85-4fe1-a353-d521f0df730f.md"
)
$wsDeDe.Range("J3").Value = "7e314904-3a85-4fe1-a353-d521f0df730f.a47c97363969567fabfb364e0821a2933e6f6093.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-28 12:24:57"

# ---------------------------------------------------------------------------
# 4. Widen columns that now hold long hyperlink / filename text so everything
#    stays readable (Overview E:F, zh-cn/de-de column C, and the newly
#    populated I:J "Latest Target File" / "Latest Handback File" columns).
# ---------------------------------------------------------------------------
$wsOverview.Range("E1:F1").ColumnWidth = 29.9777047293527

$wsZhCn.Range("C1").ColumnWidth = 29.9777047293527
$wsZhCn.Range("I1:J1").ColumnWidth = 40

$wsDeDe.Range("C1").ColumnWidth = 29.9777047293527
$wsDeDe.Range("I1:J1").ColumnWidth = 40
